$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.391.82"
$ws.Range("E2").Value = "  -6.17%  "
$ws.Range("D3").Value = "2.517.16"
$ws.Range("E3").Value = "  -3.46%  "
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").Value = "'296.14"
$ws.Range("E5").Value = "  -3.91%  "
$ws.Range("D6").Value = "'94.26"
$ws.Range("E6").Value = "  -5.95%  "
$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  -5.09%  "
$ws.Range("E8").Value = "  -0.13%  "
$ws.Range("D9").Value = "'0.550"
$ws.Range("E9").Value = "  -5.34%  "
$ws.Range("D10").Value = "'36.45"
$ws.Range("E10").Value = "  -7.70%  "
$ws.Range("D11").Value = "'0.0802"
$ws.Range("E11").Value = "  -5.04%  "
$ws.Range("D12").Value = "'7.67"
$ws.Range("E12").Value = "  -5.94%  "
$ws.Range("E13").Value = "  +0.92%  "
$ws.Range("D14").Value = "2.900.56"
$ws.Range("E14").Value = "  -3.67%  "
$ws.Range("D15").Value = "2.519.91"
$ws.Range("E15").Value = "  -3.20%  "
$ws.Range("D16").Value = "'0.869"
$ws.Range("E16").Value = "  -5.63%  "
$ws.Range("D17").Value = "'14.05"
$ws.Range("E17").Value = "  -5.85%  "
$ws.Range("D18").Value = "43.355.95"
$ws.Range("E18").Value = "  -6.70%  "
$ws.Range("D19").Value = "0.0₃0964"
$ws.Range("E19").Value = "  -4.90%  "
$ws.Range("D20").Value = "'6.54"
$ws.Range("E20").Value = "  -2.95%  "
$ws.Range("D21").Value = "'12.30"
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").Value = "'72.10"
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("D23").Value = "'259.48"
$ws.Range("E23").Value = "  -4.88%  "
$ws.Range("D24").Value = "'2.89"
$ws.Range("E24").Value = "  -4.83%  "
$ws.Range("D25").Value = "'2.14"
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").Value = "'28.74"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").Value = "'10.01"
$ws.Range("E28").Value = "  -5.51%  "
$ws.Range("D29").Value = "'2.23"
$ws.Range("E29").Value = "  -1.90%  "
$ws.Range("D30").Value = "'36.91"
$ws.Range("E30").Value = "  -5.26%  "
$ws.Range("D31").Value = "'6.04"
$ws.Range("E31").Value = "  -4.68%  "
$ws.Range("D32").Value = "'3.45"
$ws.Range("E32").Value = "  -5.57%  "
$ws.Range("B33").Value = "WEMIXToken"
$ws.Range("C33").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D33").Value = "'2.77"
$ws.Range("E33").Value = "  -3.32%  "
$ws.Range("B34").Value = "Monero"
$ws.Range("C34").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D34").Value = "'149.99"
$ws.Range("E34").Value = "  -0.67%  "
$ws.Range("D35").Value = "'2.14"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").Value = "'0.0797"
$ws.Range("E36").Value = "  -4.98%  "
$ws.Range("D37").Value = "'0.114"
$ws.Range("E37").Value = "  -6.02%  "
$ws.Range("E38").Value = "  -3.71%  "
$ws.Range("D39").Value = "'23.55"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("D40").Value = "'16.10"
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("D41").Value = "'3.47"
$ws.Range("E41").Value = "  -4.85%  "
$ws.Range("D42").Value = "'0.0308"
$ws.Range("E42").Value = "  -7.13%  "
$ws.Range("D43").Value = "'3.79"
$ws.Range("E43").Value = "  -6.98%  "
$ws.Range("D44").Value = "2.009.69"
$ws.Range("E44").Value = "  -5.34%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("D46").Value = "'85.62"
$ws.Range("E46").Value = "  -8.35%  "
$ws.Range("D47").Value = "'1.63"
$ws.Range("E47").Value = "  +4.55%  "
$ws.Range("D48").Value = "'8.90"
$ws.Range("E48").Value = "  -6.44%  "
$ws.Range("D49").Value = "2.775.29"
$ws.Range("E49").Value = "  -3.19%  "
$ws.Range("D50").Value = "'103.07"
$ws.Range("E50").Value = "  -5.55%  "
$ws.Range("E51").Value = "  -6.33%  "
